$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Ligand detection rate, Ligand average/total expression value (F, G, H)
# and Receptor/Edge derived-specificity metrics (M-T) with refreshed TPM-based numbers.
# Row 2
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.140313
$ws.Range("H2").Value2 = 0.420939
$ws.Range("M2").Value2 = 6.2878035
$ws.Range("N2").Value2 = 12.575607
$ws.Range("O2").Value2 = 0.4295541170219724
$ws.Range("P2").Value2 = 0.3581429706649357
$ws.Range("Q2").Value2 = 0.8822605724954999
$ws.Range("R2").Value2 = 5.293563434973
$ws.Range("S2").Value2 = 0.4295541170219724
$ws.Range("T2").Value2 = 0.3581429706649357

# Row 3
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.140313
$ws.Range("H3").Value2 = 0.420939
$ws.Range("M3").Value2 = 1.688011333333333
$ws.Range("N3").Value2 = 5.064033999999999
$ws.Range("O3").Value2 = 0.1153172515351477
$ws.Range("P3").Value2 = 0.1442195339205683
$ws.Range("Q3").Value2 = 0.2368499342139999
$ws.Range("R3").Value2 = 2.131649407926
$ws.Range("S3").Value2 = 0.1153172515351477
$ws.Range("T3").Value2 = 0.1442195339205683

# Row 4
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.140313
$ws.Range("H4").Value2 = 0.420939
$ws.Range("M4").Value2 = 1.53237
$ws.Range("N4").Value2 = 4.59711
$ws.Range("O4").Value2 = 0.1046845440225605
$ws.Range("P4").Value2 = 0.1309219214526569
$ws.Range("Q4").Value2 = 0.21501143181
$ws.Range("R4").Value2 = 1.93510288629
$ws.Range("S4").Value2 = 0.1046845440225605
$ws.Range("T4").Value2 = 0.1309219214526569

# Row 5
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.140313
$ws.Range("H5").Value2 = 0.420939
$ws.Range("M5").Value2 = 2.512758
$ws.Range("N5").Value2 = 5.025516
$ws.Range("O5").Value2 = 0.1716601900774885
$ws.Range("P5").Value2 = 0.1431225728797158
$ws.Range("Q5").Value2 = 0.352572613254
$ws.Range("R5").Value2 = 2.115435679524
$ws.Range("S5").Value2 = 0.1716601900774885
$ws.Range("T5").Value2 = 0.1431225728797158

# Row 6
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.140313
$ws.Range("H6").Value2 = 0.420939
$ws.Range("M6").Value2 = 1.749515666666667
$ws.Range("N6").Value2 = 5.248547
$ws.Range("O6").Value2 = 0.1195189476597205
$ws.Range("P6").Value2 = 0.1494743127909878
$ws.Range("Q6").Value2 = 0.245479791737
$ws.Range("R6").Value2 = 2.209318125633
$ws.Range("S6").Value2 = 0.1195189476597205
$ws.Range("T6").Value2 = 0.1494743127909878

# Row 7
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.140313
$ws.Range("H7").Value2 = 0.420939
$ws.Range("M7").Value2 = 0.867519
$ws.Range("N7").Value2 = 2.602557
$ws.Range("O7").Value2 = 0.05926494968311025
$ws.Range("P7").Value2 = 0.0741186882911356
$ws.Range("Q7").Value2 = 0.121724193447
$ws.Range("R7").Value2 = 1.095517741023
$ws.Range("S7").Value2 = 0.05926494968311025
$ws.Range("T7").Value2 = 0.0741186882911356

